# Fruta / hortaliza, semanal
# Insert a new weekly sample row at row 9 (pushing the existing rows 9-74
# down to 10-75) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 9..74 down to 10..75, keeping formatting intact.
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the newly inserted row 9 with the new weekly observation.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C9").Value = 'Ñuble'
$ws.Range("D9").Value = 44649
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112021
$ws.Range("G9").Value = 'Ají'
$ws.Range("H9").Value = 'Americana (o)'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 8500
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 8750
$ws.Range("N9").Value = '$/caja 15 kilos'
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 583
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = 'Hortaliza'
